# PHILIBERT_JournalDeBord.xlsx - add two new journal entries (rows 29 & 30)
# describing the Raspberry Pi server installation/configuration work.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Feuil1")

# --- Row 29 : "Réalisation" / Installation du Raspberry Pi ---------------
# Fill the "Commentaires" and "Description" cells first so the new shared
# strings end up in the same order as in the reference workbook.
$ws.Cells.Item(29, 5).Value2 = "J'installe les différents composants nécessaire au bon fonctionnement du serveur"
$ws.Cells.Item(30, 2).Value2 = "Documentation de l'installation et de la configuration du Raspberry Pi"
$ws.Cells.Item(29, 2).Value2 = "Installation du Raspberry Pi et configuration de fail2ban, ufw, ssh, nodejs et npm"

$ws.Cells.Item(29, 1).WrapText = $true
$ws.Cells.Item(29, 1).Value2 = "Réalisation"
$ws.Cells.Item(29, 3).Value2 = 2
$ws.Cells.Item(29, 4).Value2 = 43514

# --- Row 30 : "Documentation" / Documentation du Raspberry Pi ------------
$ws.Cells.Item(30, 1).WrapText = $true
$ws.Cells.Item(30, 1).Value2 = "Documentation"
$ws.Cells.Item(30, 3).Value2 = 1
$ws.Cells.Item(30, 4).Value2 = 43514

# Match the active selection shown in the source workbook after the edit.
$ws.Range("B30").Select()
